# Plan de pruebas v1.0.1
# Adds two new test cases (CP-005 Checkboxes, CP-006 Context menu) to rows 7 and 8,
# adds a new "En proceso" conditional-formatting state, and updates the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7: CP-005 Checkboxes -------------------------------------------------
$ws.Cells.Item(7, 2).Value = "CP-005 Checkboxes"
$ws.Cells.Item(7, 3).Value = "1. Ingresar en la url ""https://the-internet.herokuapp.com/basic_auth""`n2. Dar click a boton ""Checkboxes""`n3. Validar los 2 checkbox"
$ws.Cells.Item(7, 4).Value = "Se puede marcar y desmarcar los 2 checkbox"
$ws.Cells.Item(7, 6).Value = "En proceso"
$ws.Rows.Item(7).RowHeight = 67.5

# --- Row 8: CP-006 Context menu -----------------------------------------------
$ws.Cells.Item(8, 2).Value = "CP-006 Contxt menu"
$ws.Cells.Item(8, 3).Value = "1. Ingresar en la url ""https://the-internet.herokuapp.com/basic_auth""`n2. Dar click a boton ""Context menu""`n3. Dar click derecho a la caja`n4. Validar el despliegue del menu de contexto"
$ws.Cells.Item(8, 4).Value = "Al dar click derecho a la caja se despliega un menu contextual"
$ws.Cells.Item(8, 6).Value = "En proceso"
$ws.Rows.Item(8).RowHeight = 90

# --- Id caso de prueba column, filled in after the rest ----------------------
$ws.Cells.Item(7, 5).Value = "CP-005"
$ws.Cells.Item(8, 5).Value = "CP-006"

# --- New conditional formatting rule for "En proceso" -------------------------
# xlTextString = 9, xlContains = 0; the search text is passed positionally as
# the "String" argument (5th parameter) so the generated formula matches
# Excel's own SEARCH()-based "text contains" rule.
$range = $ws.Range("F3:F38")
$fc = $range.FormatConditions.Add(9, 0, [Type]::Missing, [Type]::Missing, "En proceso")
$fc.Interior.Color = 15773696
$fc.SetFirstPriority()

# --- Update the selected cell/range in the sheet view -------------------------
$ws.Range("B1:G1").Select()

$wb.Save()
